$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(688, 1).NumberFormat = "@"
$ws.Cells.Item(688, 1).Value = "2026/01/24"
$ws.Cells.Item(688, 1).Style = "Normal"
$ws.Cells.Item(688, 2).Value = "土"
$ws.Cells.Item(688, 3).Value = 16
$ws.Cells.Item(688, 4).Value = 14

$ws.Cells.Item(689, 1).NumberFormat = "@"
$ws.Cells.Item(689, 1).Value = "2026/12/29"
$ws.Cells.Item(689, 1).Style = "Normal"
$ws.Cells.Item(689, 2).Value = "火"
$ws.Cells.Item(689, 3).Value = 13
$ws.Cells.Item(689, 4).Value = 9

$ws.Cells.Item(690, 1).NumberFormat = "@"
$ws.Cells.Item(690, 1).Value = "2026/12/29"
$ws.Cells.Item(690, 1).Style = "Normal"
$ws.Cells.Item(690, 2).Value = "火"
$ws.Cells.Item(690, 3).Value = 16
$ws.Cells.Item(690, 4).Value = 10

$ws.Cells.Item(691, 1).NumberFormat = "@"
$ws.Cells.Item(691, 1).Value = "2026/12/29"
$ws.Cells.Item(691, 1).Style = "Normal"
$ws.Cells.Item(691, 2).Value = "火"
$ws.Cells.Item(691, 3).Value = 19
$ws.Cells.Item(691, 4).Value = 10

$ws.Cells.Item(692, 1).NumberFormat = "@"
$ws.Cells.Item(692, 1).Value = "2026/12/29"
$ws.Cells.Item(692, 1).Style = "Normal"
$ws.Cells.Item(692, 2).Value = "火"
$ws.Cells.Item(692, 3).Value = 23
$ws.Cells.Item(692, 4).Value = 10

$ws.Cells.Item(693, 1).NumberFormat = "@"
$ws.Cells.Item(693, 1).Value = "2026/12/30"
$ws.Cells.Item(693, 1).Style = "Normal"
$ws.Cells.Item(693, 2).Value = "水"
$ws.Cells.Item(693, 3).Value = 2
$ws.Cells.Item(693, 4).Value = 10

$ws.Cells.Item(694, 1).NumberFormat = "@"
$ws.Cells.Item(694, 1).Value = "2026/12/30"
$ws.Cells.Item(694, 1).Style = "Normal"
$ws.Cells.Item(694, 2).Value = "水"
$ws.Cells.Item(694, 3).Value = 5
$ws.Cells.Item(694, 4).Value = 10

$ws.Cells.Item(695, 1).NumberFormat = "@"
$ws.Cells.Item(695, 1).Value = "2026/12/30"
$ws.Cells.Item(695, 1).Style = "Normal"
$ws.Cells.Item(695, 2).Value = "水"
$ws.Cells.Item(695, 3).Value = 8
$ws.Cells.Item(695, 4).Value = 10

$ws.Cells.Item(696, 1).NumberFormat = "@"
$ws.Cells.Item(696, 1).Value = "2026/12/30"
$ws.Cells.Item(696, 1).Style = "Normal"
$ws.Cells.Item(696, 2).Value = "水"
$ws.Cells.Item(696, 3).Value = 13
$ws.Cells.Item(696, 4).Value = 10

$ws.Cells.Item(697, 1).NumberFormat = "@"
$ws.Cells.Item(697, 1).Value = "2026/12/30"
$ws.Cells.Item(697, 1).Style = "Normal"
$ws.Cells.Item(697, 2).Value = "水"
$ws.Cells.Item(697, 3).Value = 16
$ws.Cells.Item(697, 4).Value = 9

$ws.Cells.Item(698, 1).NumberFormat = "@"
$ws.Cells.Item(698, 1).Value = "2026/12/30"
$ws.Cells.Item(698, 1).Style = "Normal"
$ws.Cells.Item(698, 2).Value = "水"
$ws.Cells.Item(698, 3).Value = 22
$ws.Cells.Item(698, 4).Value = 10

$ws.Cells.Item(699, 1).NumberFormat = "@"
$ws.Cells.Item(699, 1).Value = "2026/12/31"
$ws.Cells.Item(699, 1).Style = "Normal"
$ws.Cells.Item(699, 2).Value = "木"
$ws.Cells.Item(699, 3).Value = 2
$ws.Cells.Item(699, 4).Value = 11

$ws.Cells.Item(700, 1).NumberFormat = "@"
$ws.Cells.Item(700, 1).Value = "2026/12/31"
$ws.Cells.Item(700, 1).Style = "Normal"
$ws.Cells.Item(700, 2).Value = "木"
$ws.Cells.Item(700, 3).Value = 6
$ws.Cells.Item(700, 4).Value = 13

$ws.Cells.Item(701, 1).NumberFormat = "@"
$ws.Cells.Item(701, 1).Value = "2026/12/31"
$ws.Cells.Item(701, 1).Style = "Normal"
$ws.Cells.Item(701, 2).Value = "木"
$ws.Cells.Item(701, 3).Value = 9
$ws.Cells.Item(701, 4).Value = 13

$ws.Cells.Item(702, 1).NumberFormat = "@"
$ws.Cells.Item(702, 1).Value = "2026/12/31"
$ws.Cells.Item(702, 1).Style = "Normal"
$ws.Cells.Item(702, 2).Value = "木"
$ws.Cells.Item(702, 3).Value = 12
$ws.Cells.Item(702, 4).Value = 34

$ws.Cells.Item(703, 1).NumberFormat = "@"
$ws.Cells.Item(703, 1).Value = "2026/12/31"
$ws.Cells.Item(703, 1).Style = "Normal"
$ws.Cells.Item(703, 2).Value = "木"
$ws.Cells.Item(703, 3).Value = 14
$ws.Cells.Item(703, 4).Value = 15

$ws.Cells.Item(704, 1).NumberFormat = "@"
$ws.Cells.Item(704, 1).Value = "2026/12/31"
$ws.Cells.Item(704, 1).Style = "Normal"
$ws.Cells.Item(704, 2).Value = "木"
$ws.Cells.Item(704, 3).Value = 22
$ws.Cells.Item(704, 4).Value = 12

$ws.Cells.Item(705, 1).NumberFormat = "@"
$ws.Cells.Item(705, 1).Value = "2027/01/01"
$ws.Cells.Item(705, 1).Style = "Normal"
$ws.Cells.Item(705, 2).Value = "金"
$ws.Cells.Item(705, 3).Value = 2
$ws.Cells.Item(705, 4).Value = 13

$ws.Cells.Item(706, 1).NumberFormat = "@"
$ws.Cells.Item(706, 1).Value = "2027/01/01"
$ws.Cells.Item(706, 1).Style = "Normal"
$ws.Cells.Item(706, 2).Value = "金"
$ws.Cells.Item(706, 3).Value = 5
$ws.Cells.Item(706, 4).Value = 12

$ws.Cells.Item(707, 1).NumberFormat = "@"
$ws.Cells.Item(707, 1).Value = "2027/01/01"
$ws.Cells.Item(707, 1).Style = "Normal"
$ws.Cells.Item(707, 2).Value = "金"
$ws.Cells.Item(707, 3).Value = 13
$ws.Cells.Item(707, 4).Value = 14

$ws.Cells.Item(708, 1).NumberFormat = "@"
$ws.Cells.Item(708, 1).Value = "2027/01/01"
$ws.Cells.Item(708, 1).Style = "Normal"
$ws.Cells.Item(708, 2).Value = "金"
$ws.Cells.Item(708, 3).Value = 16
$ws.Cells.Item(708, 4).Value = 11

$ws.Cells.Item(709, 1).NumberFormat = "@"
$ws.Cells.Item(709, 1).Value = "2027/01/01"
$ws.Cells.Item(709, 1).Style = "Normal"
$ws.Cells.Item(709, 2).Value = "金"
$ws.Cells.Item(709, 3).Value = 19
$ws.Cells.Item(709, 4).Value = 13

$ws.Cells.Item(710, 1).NumberFormat = "@"
$ws.Cells.Item(710, 1).Value = "2027/01/02"
$ws.Cells.Item(710, 1).Style = "Normal"
$ws.Cells.Item(710, 2).Value = "土"
$ws.Cells.Item(710, 3).Value = 1
$ws.Cells.Item(710, 4).Value = 12

$ws.Cells.Item(711, 1).NumberFormat = "@"
$ws.Cells.Item(711, 1).Value = "2027/01/02"
$ws.Cells.Item(711, 1).Style = "Normal"
$ws.Cells.Item(711, 2).Value = "土"
$ws.Cells.Item(711, 3).Value = 5
$ws.Cells.Item(711, 4).Value = 12

$ws.Cells.Item(712, 1).NumberFormat = "@"
$ws.Cells.Item(712, 1).Value = "2027/01/02"
$ws.Cells.Item(712, 1).Style = "Normal"
$ws.Cells.Item(712, 2).Value = "土"
$ws.Cells.Item(712, 3).Value = 8
$ws.Cells.Item(712, 4).Value = 13

$ws.Cells.Item(713, 1).NumberFormat = "@"
$ws.Cells.Item(713, 1).Value = "2027/01/02"
$ws.Cells.Item(713, 1).Style = "Normal"
$ws.Cells.Item(713, 2).Value = "土"
$ws.Cells.Item(713, 3).Value = 13
$ws.Cells.Item(713, 4).Value = 16

$ws.Cells.Item(714, 1).NumberFormat = "@"
$ws.Cells.Item(714, 1).Value = "2027/01/02"
$ws.Cells.Item(714, 1).Style = "Normal"
$ws.Cells.Item(714, 2).Value = "土"
$ws.Cells.Item(714, 3).Value = 16
$ws.Cells.Item(714, 4).Value = 19

$ws.Cells.Item(715, 1).NumberFormat = "@"
$ws.Cells.Item(715, 1).Value = "2027/01/02"
$ws.Cells.Item(715, 1).Style = "Normal"
$ws.Cells.Item(715, 2).Value = "土"
$ws.Cells.Item(715, 3).Value = 19
$ws.Cells.Item(715, 4).Value = 21

$ws.Cells.Item(716, 1).NumberFormat = "@"
$ws.Cells.Item(716, 1).Value = "2027/01/02"
$ws.Cells.Item(716, 1).Style = "Normal"
$ws.Cells.Item(716, 2).Value = "土"
$ws.Cells.Item(716, 3).Value = 22
$ws.Cells.Item(716, 4).Value = 22

$ws.Cells.Item(717, 1).NumberFormat = "@"
$ws.Cells.Item(717, 1).Value = "2027/01/03"
$ws.Cells.Item(717, 1).Style = "Normal"
$ws.Cells.Item(717, 2).Value = "日"
$ws.Cells.Item(717, 3).Value = 1
$ws.Cells.Item(717, 4).Value = 23

$ws.Cells.Item(718, 1).NumberFormat = "@"
$ws.Cells.Item(718, 1).Value = "2027/01/03"
$ws.Cells.Item(718, 1).Style = "Normal"
$ws.Cells.Item(718, 2).Value = "日"
$ws.Cells.Item(718, 3).Value = 4
$ws.Cells.Item(718, 4).Value = 26

$ws.Cells.Item(719, 1).NumberFormat = "@"
$ws.Cells.Item(719, 1).Value = "2027/01/03"
$ws.Cells.Item(719, 1).Style = "Normal"
$ws.Cells.Item(719, 2).Value = "日"
$ws.Cells.Item(719, 3).Value = 7
$ws.Cells.Item(719, 4).Value = 23

$ws.Cells.Item(720, 1).NumberFormat = "@"
$ws.Cells.Item(720, 1).Value = "2027/01/03"
$ws.Cells.Item(720, 1).Style = "Normal"
$ws.Cells.Item(720, 2).Value = "日"
$ws.Cells.Item(720, 3).Value = 13
$ws.Cells.Item(720, 4).Value = 23

$ws.Cells.Item(721, 1).NumberFormat = "@"
$ws.Cells.Item(721, 1).Value = "2027/01/03"
$ws.Cells.Item(721, 1).Style = "Normal"
$ws.Cells.Item(721, 2).Value = "日"
$ws.Cells.Item(721, 3).Value = 16
$ws.Cells.Item(721, 4).Value = 24

$ws.Cells.Item(722, 1).NumberFormat = "@"
$ws.Cells.Item(722, 1).Value = "2027/01/03"
$ws.Cells.Item(722, 1).Style = "Normal"
$ws.Cells.Item(722, 2).Value = "日"
$ws.Cells.Item(722, 3).Value = 19
$ws.Cells.Item(722, 4).Value = 26

$ws.Cells.Item(723, 1).NumberFormat = "@"
$ws.Cells.Item(723, 1).Value = "2027/01/03"
$ws.Cells.Item(723, 1).Style = "Normal"
$ws.Cells.Item(723, 2).Value = "日"
$ws.Cells.Item(723, 3).Value = 22
$ws.Cells.Item(723, 4).Value = 21

$ws.Cells.Item(724, 1).NumberFormat = "@"
$ws.Cells.Item(724, 1).Value = "2027/01/04"
$ws.Cells.Item(724, 1).Style = "Normal"
$ws.Cells.Item(724, 2).Value = "月"
$ws.Cells.Item(724, 3).Value = 2
$ws.Cells.Item(724, 4).Value = 19

$ws.Cells.Item(725, 1).NumberFormat = "@"
$ws.Cells.Item(725, 1).Value = "2027/01/04"
$ws.Cells.Item(725, 1).Style = "Normal"
$ws.Cells.Item(725, 2).Value = "月"
$ws.Cells.Item(725, 3).Value = 4
$ws.Cells.Item(725, 4).Value = 18

$ws.Cells.Item(726, 1).NumberFormat = "@"
$ws.Cells.Item(726, 1).Value = "2027/01/04"
$ws.Cells.Item(726, 1).Style = "Normal"
$ws.Cells.Item(726, 2).Value = "月"
$ws.Cells.Item(726, 3).Value = 7
$ws.Cells.Item(726, 4).Value = 19

$ws.Cells.Item(727, 1).NumberFormat = "@"
$ws.Cells.Item(727, 1).Value = "2027/01/04"
$ws.Cells.Item(727, 1).Style = "Normal"
$ws.Cells.Item(727, 2).Value = "月"
$ws.Cells.Item(727, 3).Value = 13
$ws.Cells.Item(727, 4).Value = 20

$ws.Cells.Item(728, 1).NumberFormat = "@"
$ws.Cells.Item(728, 1).Value = "2027/01/04"
$ws.Cells.Item(728, 1).Style = "Normal"
$ws.Cells.Item(728, 2).Value = "月"
$ws.Cells.Item(728, 3).Value = 22
$ws.Cells.Item(728, 4).Value = 13

$ws.Cells.Item(729, 1).NumberFormat = "@"
$ws.Cells.Item(729, 1).Value = "2027/01/05"
$ws.Cells.Item(729, 1).Style = "Normal"
$ws.Cells.Item(729, 2).Value = "火"
$ws.Cells.Item(729, 3).Value = 1
$ws.Cells.Item(729, 4).Value = 13

$ws.Cells.Item(730, 1).NumberFormat = "@"
$ws.Cells.Item(730, 1).Value = "2027/01/05"
$ws.Cells.Item(730, 1).Style = "Normal"
$ws.Cells.Item(730, 2).Value = "火"
$ws.Cells.Item(730, 3).Value = 7
$ws.Cells.Item(730, 4).Value = 14
